# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets
# to reflect newly generated output (gh-pages regeneration at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1 / sheetId 1)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 188
$wsExhibit.Range("F5").Value = 5123
$wsExhibit.Range("F6").Value = 26
$wsExhibit.Range("F9").Value = 569
$wsExhibit.Range("F10").Value = 524
$wsExhibit.Range("F13").Value = 1427
$wsExhibit.Range("F14").Value = 3915
$wsExhibit.Range("F15").Value = 422
$wsExhibit.Range("F16").Value = 154
$wsExhibit.Range("F17").Value = 137
$wsExhibit.Range("F18").Value = 90
$wsExhibit.Range("F19").Value = 3074
$wsExhibit.Range("F20").Value = 147
$wsExhibit.Range("F21").Value = 1040
$wsExhibit.Range("F22").Value = 95
$wsExhibit.Range("F23").Value = 40
$wsExhibit.Range("F24").Value = 183
$wsExhibit.Range("F25").Value = 79
$wsExhibit.Range("F30").Value = 6
$wsExhibit.Range("F32").Value = 7
$wsExhibit.Range("F33").Value = 3

# Sheet "全部类型" (index 4 / sheetId 4)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 188
$wsAll.Range("F6").Value = 5123
$wsAll.Range("F7").Value = 26
$wsAll.Range("F10").Value = 569
$wsAll.Range("F11").Value = 524
$wsAll.Range("F14").Value = 1427
$wsAll.Range("F15").Value = 3915
$wsAll.Range("F16").Value = 422
$wsAll.Range("F17").Value = 154
$wsAll.Range("F18").Value = 137
$wsAll.Range("F19").Value = 90
$wsAll.Range("F20").Value = 3074
$wsAll.Range("F21").Value = 147
$wsAll.Range("F22").Value = 1040
$wsAll.Range("F23").Value = 95
$wsAll.Range("F24").Value = 40
$wsAll.Range("F25").Value = 183
$wsAll.Range("F26").Value = 79
$wsAll.Range("F31").Value = 6
$wsAll.Range("F33").Value = 7
$wsAll.Range("F34").Value = 3
